$wb = $excel.ActiveWorkbook

# Add the new sheet after the existing InstallGames sheet, matching
# the target workbook's sheet order.
$installGamesSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $installGamesSheet)
$ws.Name = "Sheet1"

$data = @(
    @('Column1', 'Column2'),
    @('Nomads of Driftland', 'https://store.steampowered.com/app/1400520/Nomads_of_Driftland/?snr=1_1510_4__103'),
    @('Slapshot: Rebound', 'https://store.steampowered.com/app/1173370/Slapshot_Rebound/?snr=1_1510_4__103'),
    @('Leaf Blower Revolution - Idle Game', 'https://store.steampowered.com/app/1468260/Leaf_Blower_Revolution__Idle_Game/?snr=1_1510_4__103'),
    @('Summerland', 'https://store.steampowered.com/app/1444110/Summerland/?snr=1_1510_4__103'),
    @('Ruff Night At The Gallery', 'https://store.steampowered.com/app/1463470/Ruff_Night_At_The_Gallery/?snr=1_1510_4__103'),
    @('Rubber Bandits: Christmas Prologue', 'https://store.steampowered.com/app/1469250/Rubber_Bandits_Christmas_Prologue/?snr=1_1510_4__103'),
    @('The Life and Suffering of Sir Brante — Chapter 1&2', 'https://store.steampowered.com/app/1460090/The_Life_and_Suffering_of_Sir_Brante__Chapter_12/?snr=1_1510_4__103'),
    @('AOD: Art Of Defense', 'https://store.steampowered.com/app/1367760/AOD_Art_Of_Defense/?snr=1_1510_4__103'),
    @('Mrs. Santa''s Gift Hunt', 'https://store.steampowered.com/app/1473660/Mrs_Santas_Gift_Hunt/?snr=1_1510_4__103'),
    @('Ninja Hanrei', 'https://store.steampowered.com/app/1384620/Ninja_Hanrei/?snr=1_1510_4__103'),
    @('*NEW* SCUFFED EPIC BHOP SIMULATOR 2023 (POG CHAMP)', 'https://store.steampowered.com/app/1465860/NEW_SCUFFED_EPIC_BHOP_SIMULATOR_2023_POG_CHAMP/?snr=1_1510_4__103'),
    @('DreamWatcher', 'https://store.steampowered.com/app/1444750/DreamWatcher/?snr=1_1510_4__103'),
    @('From the Shadows', 'https://store.steampowered.com/app/1427210/From_the_Shadows/?snr=1_1510_4__103'),
    @('KAKU: Ancient Seal (Alpha)', 'https://store.steampowered.com/app/1394690/KAKU_Ancient_Seal_Alpha/?snr=1_1510_4__103'),
    @('Night Reverie: Prologue', 'https://store.steampowered.com/app/1453600/Night_Reverie_Prologue/?snr=1_1510_4__103')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Output ("rows written: " + $data.Length)
